# Daily attendance processing: re-order the comma-separated "Recorded By"
# (column G) editor list on each row according to a fixed editor-priority
# ranking (derived from the canonical reordering observed for this sheet):
#   admin@admin.com < System < backup@backdoor.com < dnasr281@gmail.com < system
# Sort is stable, so unrecognised names keep their relative order at the end.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-NamePriority($name) {
  if ($name.Equals("admin@admin.com")) { return 0 }
  if ($name.Equals("System")) { return 1 }
  if ($name.Equals("backup@backdoor.com")) { return 2 }
  if ($name.Equals("dnasr281@gmail.com")) { return 3 }
  if ($name.Equals("system")) { return 4 }
  return 99
}

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
  $cell = $ws.Cells.Item($row, 7)
  $val = $cell.Value2
  if ($null -eq $val) { continue }
  if (-not ($val -is [string])) { continue }
  if ($val -eq "") { continue }
  if ($val.IndexOf(", ") -lt 0) { continue }

  $parts = $val -split ", "

  # Stable sort: decorate with (priority, original-index), sort, extract.
  $decorated = @()
  for ($i = 0; $i -lt $parts.Count; $i++) {
    $p = $parts[$i]
    $pr = Get-NamePriority $p
    $decorated += , @($pr, $i, $p)
  }
  $sortedDecorated = $decorated | Sort-Object { $_[0] }, { $_[1] }
  $newParts = @()
  foreach ($d in $sortedDecorated) { $newParts += $d[2] }
  $newVal = [string]::Join(", ", $newParts)

  if (-not $newVal.Equals($val)) {
    $cell.Value = $newVal
  }
}
